# Insert a new data row at row 81 (pushing the existing rows 81-100 down
# to 82-101) and populate it with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 81..100 down to 82..101, keeping all their existing data/format.
$ws.Rows.Item(81).Insert()

# Populate the newly inserted row 81 with the new record.
$ws.Cells.Item(81, 1).Value = 8
$ws.Cells.Item(81, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(81, 3).Value = "Coquimbo"
$ws.Cells.Item(81, 4).Value = Get-Date -Year 2022 -Month 3 -Day 21 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(81, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(81, 5).Value = 4
$ws.Cells.Item(81, 6).Value = "Fruta"
$ws.Cells.Item(81, 7).Value = 100109
$ws.Cells.Item(81, 8).Value = "Uva"
$ws.Cells.Item(81, 9).Value = 100109001
$ws.Cells.Item(81, 10).Value = "Uva"
$ws.Cells.Item(81, 11).Value = "Red Globe"
$ws.Cells.Item(81, 12).Value = "Primera"
$ws.Cells.Item(81, 13).Value = 400
$ws.Cells.Item(81, 14).Value = 9500
$ws.Cells.Item(81, 15).Value = 10000
$ws.Cells.Item(81, 16).Value = 9750
$ws.Cells.Item(81, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(81, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(81, 19).Value = 542
$ws.Cells.Item(81, 20).Value = 18
